$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41..84 down to 42..85.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new record's data.
$ws.Range("A41").Value = 10
$ws.Range("B41").Value = "Vega Modelo de Temuco"
$ws.Range("C41").Value = "La Araucanía"
$ws.Range("D41").Value = 44566
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = 100112022
$ws.Range("G41").Value = "Arveja Verde"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 20
$ws.Range("K41").Value = 26000
$ws.Range("L41").Value = 26000
$ws.Range("M41").Value = 26000
$ws.Range("N41").Value = "`$/saco 25 kilos"
$ws.Range("O41").Value = "Región de La Araucanía"
$ws.Range("P41").Value = 1040
$ws.Range("Q41").Value = 25
$ws.Range("R41").Value = "Hortaliza"
